# CAMPO_NOVO.xlsx update
# - Rename "Paineis DARQ" -> "PAINEIS DARQ"
# - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Delete the "Desarquivamentos Pendentes" sheet

$wb = $excel.ActiveWorkbook

# Avoid any confirmation prompts when deleting a sheet
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Deleting the trailing sheet shifts the active tab; restore the original
# active sheet (first tab) so bookViews/tabSelected stay as before.
$wb.Worksheets.Item("PAINEIS DARQ").Activate()

$excel.DisplayAlerts = $true
